$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 ("ID" row) -- everything below shifts up by one row.
$ws.Rows.Item(2).Delete()

# Update rounded p-values that changed slightly after the shift.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "9e-48"   # index_group: IR
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "9e-48"   # index_group: OR_K
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "9e-48"   # index_group: SM
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1e-18"   # grade
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3e-20"   # renal_pres: Y

# AKI: Y row values changed
$ws.Range("B8").Value = "15 (24%)"
$ws.Range("C8").Value = "39 (45%)"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.013"

# return_ed_30d: Y row values changed
$ws.Range("B9").Value = "17 (27%)"
$ws.Range("C9").Value = "30 (34%)"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.42"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3e-04"  # surv_ICU_LOS
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5e-05"  # surv_hosp_LOS

# Add new column F with header "OR" (odds ratio), styled like the header row.
$ws.Range("F1").Value = "OR"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").HorizontalAlignment = -4108

$ws.Range("F6").Value = "0.20 [0.04–0.90]"   # survived: Y
$ws.Range("F7").Value = "0.00 [0.00–NaN]"    # renal_pres: Y
$ws.Range("F8").Value = "2.60 [1.27–5.33]"   # AKI: Y
$ws.Range("F9").Value = "1.42 [0.70–2.90]"   # return_ed_30d: Y
